# Add the "purchase" sheet right after "products" and populate its header row.
$wb = $excel.ActiveWorkbook
$products = $wb.Worksheets("products")

$purchase = $wb.Worksheets.Add($null, $products)
$purchase.Name = "purchase"
$purchase.Range("A1").Value = "id"
$purchase.Range("B1").Value = "customerkey"
$purchase.Range("C1").Value = "product"
$purchase.Range("D1").Value = "quantity"
$purchase.Range("H18").Select() | Out-Null

# Add the "customers" sheet right after "purchase" and populate it.
$customers = $wb.Worksheets.Add($null, $purchase)
$customers.Name = "customers"

$customers.Range("A1").Value = "customerkey"
$customers.Range("B1").Value = "firstname"
$customers.Range("C1").Value = "lastname"

$customers.Range("B2").Value = "Ebuka"
$customers.Range("C2").Value = "Ifechukwu"
$customers.Range("A2").Formula = "=LEFT(B2,3)&LEFT(C2,3)"

$customers.Range("B3").Value = "Chidera"
$customers.Range("C3").Value = "Ifechukwu"
$customers.Range("B4").Value = "Ugonna"
$customers.Range("C4").Value = "Ifechukwu"
$customers.Range("B5").Value = "Uchenna"
$customers.Range("C5").Value = "Ifechukwu"

# One shared formula covering A3:A5 (mirrors Excel's own fill-down grouping).
$customers.Range("A3:A5").Formula = "=LEFT(B3,3)&LEFT(C3,3)"

$customers.Columns("A:B").AutoFit() | Out-Null
$customers.Range("B9").Select() | Out-Null

# Re-activate "purchase" so it ends up the selected/visible tab.
$purchase.Activate() | Out-Null
